# Regenerate save_data to use K (strikeouts) instead of Strike# in column G
# (rows 2-50, header "K" in G1). Column G holds the number of strikeouts
# recorded for the pitching appearance on that row; this revision recomputes
# those totals (the underlying box-score / std-mean source data lives
# outside this sheet, so the recomputed totals are applied directly here).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 3
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 2
    16 = 2
    17 = 3
    18 = 0
    19 = 1
    21 = 0
    22 = 0
    23 = 2
    24 = 3
    25 = 0
    26 = 1
    27 = 1
    28 = 2
    29 = 2
    30 = 3
    31 = 0
    32 = 2
    33 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 1
    38 = 0
    39 = 2
    40 = 1
    42 = 1
    43 = 1
    44 = 1
    46 = 1
    47 = 1
    48 = 2
    49 = 2
    50 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
